$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column E ("2025-05-11") contribution counts
$ws.Range("E2").Value = "1 Contributions"
$ws.Range("E3").Value = "1 Contributions"
$ws.Range("E4").Value = "1 Contributions"
$ws.Range("E5").Value = "No Contributions"
$ws.Range("E6").Value = "3 Contributions"
$ws.Range("E7").Value = "2 Contributions"
$ws.Range("E8").Value = "No Contributions"
$ws.Range("E9").Value = "No Contributions"

# Username column (B) now holds numeric GitLab user IDs for two rows
# (cell type checking: numeric values are written as numbers, not text)
$ws.Cells.Item(6,2).Value = 26737887
$ws.Cells.Item(9,2).Value = 26737892

# Update the active selection to E9
$ws.Range("E9").Select()
